{"js": "// Remove the trailing \"Ver no Jupiter...\" / \"\u00a9 2020 ...\" footer block\n// (and the blank paragraph that immediately precedes it) that followed\n// the bibliography, leaving the blank paragraph + page-break paragraph\n// that come after it untouched.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Find the boundary paragraph (end of the bibliography) and the footer\n// paragraphs by their text so the edit is robust to any surrounding\n// paragraphs that might shift indices.\nconst items = paragraphs.items;\nlet benwartIdx = -1;\nlet jupiterIdx = -1;\nlet copyrightIdx = -1;\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n  if (t.indexOf(\"BENWART\") !== -1) {\n    benwartIdx = i;\n  }\n  if (t.indexOf(\"Ver no Jupiter\") !== -1) {\n    jupiterIdx = i;\n  }\n  if (t.indexOf(\"Powered by Jekyll\") !== -1) {\n    copyrightIdx = i;\n  }\n}\n\nconst toDelete = [];\nif (jupiterIdx !== -1) {\n  toDelete.push(jupiterIdx);\n}\nif (copyrightIdx !== -1) {\n  toDelete.push(copyrightIdx);\n}\n// The blank paragraph directly after the bibliography text and before\n// the \"Ver no Jupiter\" paragraph also needs to go.\nif (benwartIdx !== -1 && jupiterIdx === benwartIdx + 2) {\n  toDelete.push(benwartIdx + 1);\n}\n\n// Delete from the end so earlier indices stay valid.\ntoDelete.sort((a, b) => b - a);\nfor (const idx of toDelete) {\n  items[idx].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" / \"\u00a9 2020 ...\" footer block\n# (and the blank paragraph that immediately precedes it) that followed\n# the bibliography, leaving the blank paragraph + page-break paragraph\n# that come after it untouched.\n$d = $word.ActiveDocument\n$paras = $d.Paragraphs\n$count = $paras.Count\n\n$benwartIdx = -1\n$jupiterIdx = -1\n$copyrightIdx = -1\n\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $paras.Item($i).Range.Text\n    if ($t -like \"*BENWART*\") { $benwartIdx = $i }\n    if ($t -like \"*Ver no Jupiter*\") { $jupiterIdx = $i }\n    if ($t -like \"*Powered by Jekyll*\") { $copyrightIdx = $i }\n}\n\n$toDelete = New-Object System.Collections.ArrayList\nif ($jupiterIdx -ne -1) { [void]$toDelete.Add($jupiterIdx) }\nif ($copyrightIdx -ne -1) { [void]$toDelete.Add($copyrightIdx) }\n# The blank paragraph directly after the bibliography text and before\n# the \"Ver no Jupiter\" paragraph also needs to go.\nif (($benwartIdx -ne -1) -and ($jupiterIdx -eq ($benwartIdx + 2))) {\n    [void]$toDelete.Add($benwartIdx + 1)\n}\n\n# Delete from the end so earlier indices stay valid.\n$sorted = $toDelete | Sort-Object -Descending\nforeach ($idx in $sorted) {\n    $paras.Item($idx).Range.Delete()\n}\n"}
